$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add a new data row for 2021-04-02 (serial 44288) with 3385 doses.
#
#    The sheet always keeps its very last data row highlighted with the
#    "Good" (green) cell style, while every earlier data row uses "Neutral"
#    (orange). We copy the current last row (102, "Good") down into a new
#    row 103 - Excel shifts the relative formulas for us automatically and
#    row 103 keeps the "Good" formatting (exactly what the new last row
#    needs). Row 102 then gets reformatted back to a plain "Neutral" row.
# ---------------------------------------------------------------------------

$ws.Range("A102:K102").Copy() | Out-Null
$ws.Range("A103").Insert() | Out-Null

# --- restore row 102 back to being a normal ("Neutral") row -----------------
# Copy number formats from row 101 (an ordinary "Neutral" row) onto row 102
# so the underlying style indices line up with the rest of the table instead
# of minting brand-new style entries.
$ws.Range("A101:H101").Copy() | Out-Null
$ws.Range("A102:H102").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A102").Value = 44287
$ws.Range("B102").Value = 3672
$ws.Range("C102").Formula = "=(AVERAGE(B96:B102))"
$ws.Range("D102").Formula = "=(D101-B102)"
$ws.Range("E102").Formula = "=E101+B102"
$ws.Range("F102").Formula = "=D102/C102"
$ws.Range("G102").Formula = "=A102+F102"
$ws.Range("H102").Formula = "=D102/84"

# row 102 no longer carries the "Total Doses" label / trailing blanks
$ws.Range("I102:K102").Clear()

# rows 98-101 no longer carry that same label either - only the new last
# row (103) keeps it.
$ws.Range("I98:I101").Clear()

# --- write the new last row (103) -------------------------------------------
# Row 103 already inherited the "Good" formatting and the I103 label text
# from the copy above, so we only need to update the day's figures; the
# dependent formulas recalculate on their own.
$ws.Range("A103").Value = 44288
$ws.Range("B103").Value = 3385

$ws.Range("J103:K103").Clear()

# ---------------------------------------------------------------------------
# 2. Fill in the previously-missing H column (daily rate / 84) for rows
#    91-97, matching the style already used lower down the H column.
# ---------------------------------------------------------------------------
$ws.Range("H98").Copy() | Out-Null
$ws.Range("H91:H97").PasteSpecial(-4122) | Out-Null     # xlPasteFormats
$ws.Range("H91:H97").Formula = "=D91/84"

# ---------------------------------------------------------------------------
# 3. Clean up the old trailer rows (now 104-107) and rebuild the summary
#    block below the table.
# ---------------------------------------------------------------------------
$ws.Range("A104:K107").Clear()

$ws.Range("A104").NumberFormat = "d-mmm"
$ws.Range("A105").NumberFormat = "d-mmm"
$ws.Range("A106").NumberFormat = "d-mmm"

$ws.Range("B105").Formula = "=SUM(B2:B102)"
$ws.Range("B106").Formula = "=B105-(98790+1320)"
$ws.Range("B108").Formula = "=SUM(B2:B103)"

# ---------------------------------------------------------------------------
# 4. Update the selected cell to match where the author left the cursor.
# ---------------------------------------------------------------------------
$ws.Range("B108").Select() | Out-Null
